# Updates the cryptocurrency price/volume table (refreshed data snapshot),
# including swapping the WEMIXTOKEN / Decentraland rows (45 <-> 46).
#
# Price-column (D) values are written with a leading apostrophe so Excel
# keeps them as literal text instead of re-parsing look-alike numbers
# (e.g. "55.00" -> 55, "0.4970" -> 0.497); the style is then reset to
# "Normal" so the cell is left with no explicit style index, matching the
# plain (unstyled) inline-string cells used throughout the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.118.33"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.77%  '
$ws.Range('D3').Value = "'1.996.45"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.91%  '
$ws.Range('E4').Value = '  +0.49%  '
$ws.Range('D5').Value = "'330.76"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('E6').Value = '  +0.60%  '
$ws.Range('D7').Value = "'0.4970"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.66%  '
$ws.Range('D8').Value = "'0.4190"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.89%  '
$ws.Range('D9').Value = "'55.00"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.26%  '
$ws.Range('D10').Value = "'0.08882"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.15%  '
$ws.Range('D11').Value = "'1.093"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.45%  '
$ws.Range('D12').Value = "'22.89"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.18%  '
$ws.Range('D13').Value = "'2.014.30"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.02%  '
$ws.Range('D14').Value = "'7.992"
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Value = "'6.416"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.34%  '
$ws.Range('D16').Value = "'1.015"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.60%  '
$ws.Range('D17').Value = "'92.52"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.77%  '
$ws.Range('D18').Value = "'0.00001105"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.79%  '
$ws.Range('D19').Value = "'0.06752"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.33%  '
$ws.Range('D20').Value = "'19.45"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.67%  '
$ws.Range('E21').Value = '  +0.77%  '
$ws.Range('D22').Value = "'5.965"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.93%  '
$ws.Range('D23').Value = "'29.131.42"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.82%  '
$ws.Range('D24').Value = "'11.97"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.64%  '
$ws.Range('D25').Value = "'2.293"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.76%  '
$ws.Range('D26').Value = "'2.240.74"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.12%  '
$ws.Range('D27').Value = "'20.77"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.33%  '
$ws.Range('D28').Value = "'156.94"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.47%  '
$ws.Range('D29').Value = "'6.262"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.10%  '
$ws.Range('D30').Value = "'2.246"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.65%  '
$ws.Range('E31').Value = '  -1.21%  '
$ws.Range('E32').Value = '  -2.11%  '
$ws.Range('D33').Value = "'0.09871"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.24%  '
$ws.Range('D34').Value = "'1.531"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.55%  '
$ws.Range('D35').Value = "'5.822"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.30%  '
$ws.Range('D36').Value = "'3.734"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.27%  '
$ws.Range('D37').Value = "'0.02413"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.88%  '
$ws.Range('D38').Value = "'1.312"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.18%  '
$ws.Range('D39').Value = "'9.056"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.46%  '
$ws.Range('D40').Value = "'0.06375"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.50%  '
$ws.Range('D41').Value = "'0.6464"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.98%  '
$ws.Range('D42').Value = "'11.55"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.53%  '
$ws.Range('D43').Value = "'0.1976"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.32%  '
$ws.Range('E44').Value = '  +0.63%  '
$ws.Range('B45').Value = 'WEMIXTOKEN'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').Value = "'1.356"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.17%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = "'0.6180"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.18%  '
$ws.Range('D47').Value = "'13.42"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.05%  '
$ws.Range('D48').Value = "'2.165"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.85%  '
$ws.Range('E49').Value = '  +8.39%  '
$ws.Range('D50').Value = "'3.491"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.30%  '
$ws.Range('D51').Value = "'2.172"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.44%  '
